$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force the Price/Volume columns to Text format so that
# numeric-looking values (e.g. "520.93") are stored as text, matching
# the source data which uses inline/shared strings, not numbers.
$priceVolRange = $ws.Range("D2:E51")
$priceVolRange.NumberFormat = "@"

$ws.Range("D2").Value = '57.565.84'
$ws.Range("E2").Value = '  +0.79%  '
$ws.Range("D3").Value = '3.095.86'
$ws.Range("E3").Value = '  +1.63%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").Value = '520.93'
$ws.Range("E5").Value = '  +1.43%  '
$ws.Range("D6").Value = '140.19'
$ws.Range("E6").Value = '  +0.60%  '
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("D8").Value = '3.096.42'
$ws.Range("E8").Value = '  +1.67%  '
$ws.Range("D9").Value = '0.432'
$ws.Range("E9").Value = '  +0.14%  '
$ws.Range("D10").Value = '7.23'
$ws.Range("E10").Value = '  -0.83%  '
$ws.Range("D11").Value = '0.108'
$ws.Range("E11").Value = '  +1.30%  '
$ws.Range("E12").Value = '  +2.89%  '
$ws.Range("D13").Value = '3.628.40'
$ws.Range("E13").Value = '  +1.32%  '
$ws.Range("E14").Value = '  +1.47%  '
$ws.Range("D15").Value = '26.11'
$ws.Range("E15").Value = '  +3.29%  '
$ws.Range("D16").Value = '0.0000163'
$ws.Range("E16").Value = '  +1.20%  '
$ws.Range("D17").Value = '57.665.49'
$ws.Range("E17").Value = '  +0.81%  '
$ws.Range("D18").Value = '3.103.35'
$ws.Range("E18").Value = '  +1.67%  '
$ws.Range("D19").Value = '6.08'
$ws.Range("E19").Value = '  -0.13%  '
$ws.Range("D20").Value = '12.76'
$ws.Range("E20").Value = '  -0.45%  '
$ws.Range("E21").Value = '  -0.15%  '
$ws.Range("D22").Value = '335.69'
$ws.Range("E22").Value = '  +1.40%  '
$ws.Range("D23").Value = '0.997'
$ws.Range("E23").Value = '  -0.21%  '
$ws.Range("D24").Value = '0.506'
$ws.Range("E24").Value = '  +1.77%  '
$ws.Range("D25").Value = '66.40'
$ws.Range("E25").Value = '  +1.35%  '
$ws.Range("D26").Value = '0.167'
$ws.Range("E26").Value = '  -0.21%  '
$ws.Range("D27").Value = '0.999'
$ws.Range("E27").Value = '  -0.01%  '
$ws.Range("D28").Value = '0.0₃0917'
$ws.Range("E28").Value = '  +2.14%  '
$ws.Range("D29").Value = '6.51'
$ws.Range("E29").Value = '  +3.48%  '
$ws.Range("D30").Value = '0.998'
$ws.Range("E30").Value = '  +0.00%  '
$ws.Range("D31").Value = '7.19'
$ws.Range("E31").Value = '  +1.17%  '
$ws.Range("D32").Value = '1.85'
$ws.Range("E32").Value = '  +2.44%  '
$ws.Range("D33").Value = '20.84'
$ws.Range("E33").Value = '  +0.58%  '
$ws.Range("E34").Value = '  +2.35%  '
$ws.Range("D35").Value = '153.31'
$ws.Range("E35").Value = '  +0.04%  '
$ws.Range("E36").Value = '  +4.87%  '
$ws.Range("D37").Value = '6.05'
$ws.Range("E37").Value = '  +3.45%  '
$ws.Range("D38").Value = '26.80'
$ws.Range("E38").Value = '  +0.47%  '
$ws.Range("D39").Value = '1.29'
$ws.Range("E39").Value = '  +3.06%  '
$ws.Range("D40").Value = '0.0664'
$ws.Range("E40").Value = '  -0.54%  '
$ws.Range("D41").Value = '3.140.33'
$ws.Range("E41").Value = '  +1.61%  '
$ws.Range("D42").Value = '0.682'
$ws.Range("E42").Value = '  +4.73%  '
$ws.Range("D43").Value = '36.77'
$ws.Range("E43").Value = '  +0.06%  '
$ws.Range("D44").Value = '3.87'
$ws.Range("E44").Value = '  +0.15%  '
$ws.Range("E45").Value = '  +0.15%  '
$ws.Range("E46").Value = '  +6.61%  '
$ws.Range("D47").Value = '2.275.93'
$ws.Range("E47").Value = '  +1.06%  '
$ws.Range("D48").Value = '0.0257'
$ws.Range("E48").Value = '  +1.29%  '
$ws.Range("D49").Value = '0.984'
$ws.Range("E49").Value = '  +7.30%  '
$ws.Range("D50").Value = '20.60'
$ws.Range("E50").Value = '  +2.20%  '
$ws.Range("D51").Value = '5.98'
$ws.Range("E51").Value = '  +2.51%  '

# Restore the default (General) formatting so no stray number-format
# style is left applied to the cells.
$priceVolRange.ClearFormats()
